# Weekly update: insert a new weekly record row for "Brócoli" at
# Terminal Hortofrutícola Agro Chillán, pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 162 (row 162 and everything below shifts down by one,
# so the former row 162 becomes row 163, ..., former row 193 becomes row 194).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly record.
$ws.Cells.Item(162, 1).Value  = 7
$ws.Cells.Item(162, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(162, 3).Value  = "Ñuble"
$ws.Cells.Item(162, 4).Value  = 44522
$ws.Cells.Item(162, 5).Value  = 16
$ws.Cells.Item(162, 6).Value  = 100112023
$ws.Cells.Item(162, 7).Value  = "Brócoli"
$ws.Cells.Item(162, 8).Value  = "Sin especificar"
$ws.Cells.Item(162, 9).Value  = "Primera"
$ws.Cells.Item(162, 10).Value = 300
$ws.Cells.Item(162, 11).Value = 700
$ws.Cells.Item(162, 12).Value = 800
$ws.Cells.Item(162, 13).Value = 750
$ws.Cells.Item(162, 14).Value = "`$/unidad"
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 750
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = "Hortaliza"
